$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.625.04'
$ws.Range("E2").Value = '  +1.77%  '

$ws.Range("D3").Value = '3.325.64'
$ws.Range("E3").Value = '  +2.59%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.85'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.72'
$ws.Range("E6").Value = '  +2.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +2.23%  '

$ws.Range("D9").Value = '3.318.52'
$ws.Range("E9").Value = '  +2.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.181'
$ws.Range("E10").Value = '  +6.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.579'
$ws.Range("E11").Value = '  +1.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.63'
$ws.Range("E12").Value = '  +4.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +1.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '691.72'
$ws.Range("E14").Value = '  +3.89%  '

$ws.Range("D15").Value = '3.871.97'
$ws.Range("E15").Value = '  +2.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.38'
$ws.Range("E16").Value = '  +2.11%  '

$ws.Range("D17").Value = '67.666.41'
$ws.Range("E17").Value = '  +1.93%  '

$ws.Range("E18").Value = '  +0.70%  '

$ws.Range("D19").Value = '3.321.64'
$ws.Range("E19").Value = '  +2.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.49'
$ws.Range("E20").Value = '  +2.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.01'
$ws.Range("E21").Value = '  +3.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.891'
$ws.Range("E22").Value = '  +2.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.52'
$ws.Range("E23").Value = '  +4.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.81'
$ws.Range("E24").Value = '  -0.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '101.09'
$ws.Range("E25").Value = '  +5.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.89'
$ws.Range("E26").Value = '  +1.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.67'
$ws.Range("E27").Value = '  +1.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.36'
$ws.Range("E28").Value = '  +4.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.94'
$ws.Range("E29").Value = '  +2.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.51'
$ws.Range("E30").Value = '  +3.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.00'
$ws.Range("E31").Value = '  +4.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '568.65'
$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.97'
$ws.Range("E33").Value = '  +1.59%  '

$ws.Range("E34").Value = '  +3.25%  '

$ws.Range("B35").Value = 'Dai'
$ws.Range("C35").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.21'
$ws.Range("E36").Value = '  +4.15%  '

$ws.Range("D37").Value = '3.699.82'
$ws.Range("E37").Value = '  -2.00%  '

$ws.Range("E38").Value = '  -2.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.00'
$ws.Range("E39").Value = '  +10.78%  '

$ws.Range("E40").Value = '  +3.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.14'
$ws.Range("E41").Value = '  +5.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.60'
$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("E43").Value = '  +3.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.333'
$ws.Range("E44").Value = '  +3.63%  '

$ws.Range("D45").Value = '0.0₃0667'
$ws.Range("E45").Value = '  +2.71%  '

$ws.Range("E46").Value = '  +2.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.63'
$ws.Range("E47").Value = '  +4.94%  '

$ws.Range("E48").Value = '  +2.15%  '

$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.32'
$ws.Range("E50").Value = '  +0.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '132.05'
$ws.Range("E51").Value = '  +3.12%  '
